$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.167.12'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '2.485.23'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.97'
$ws.Range("E5").Value = '  -1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.54'
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.41'
$ws.Range("E10").Value = '  +3.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0808'
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.28'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.09'
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").Value = '2.870.41'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("D16").Value = '2.488.37'
$ws.Range("E16").Value = '  -2.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.846'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '47.077.39'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.80'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.71'
$ws.Range("E22").Value = '  +14.02%  '
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '245.08'
$ws.Range("E24").Value = '  -2.96%  '
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.63'
$ws.Range("E27").Value = '  -3.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.27'
$ws.Range("E28").Value = '  +2.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.94'
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.31'
$ws.Range("E31").Value = '  -2.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.41'
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.18'
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.30'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0776'
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("E38").Value = '  -0.99%  '
$ws.Range("E39").Value = '  -2.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.52'
$ws.Range("E40").Value = '  +3.48%  '
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("E42").Value = '  -0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '118.49'
$ws.Range("E43").Value = '  -3.96%  '
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("D45").Value = '1.981.02'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("E47").Value = '  -6.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.02'
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.08'
$ws.Range("E50").Value = '  -6.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.42'
$ws.Range("E51").Value = '  +2.70%  '
